{"js": "// Merge the three runs of the \">>> your stuff ...\" paragraph into one run,\n// then add a new paragraph (\"Changes made by user of assignment 1\") and an\n// empty paragraph right after it.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text.indexOf(\"your stuff after this line\") !== -1\n);\n\nif (!target) {\n  throw new Error(\"Could not find the '>>> your stuff after this line >>>' paragraph\");\n}\n\n// Rewrite the paragraph's text as a single run (collapses the 3-run split\n// caused by the grammar-check proofErr markers around \"your\").\ntarget.insertText(\">>>  your stuff after this line >>>\", \"Replace\");\n\n// Insert the new paragraph with the commit text right after it...\nconst newPara = target.insertParagraph(\"Changes made by user of assignment 1\", \"After\");\n\n// ...followed by a blank paragraph (no text -> no empty run is created).\nnewPara.insertParagraph(undefined, \"After\");\n\nawait context.sync();\n", "ps1": "# Merge the three runs of the \">>> your stuff ...\" paragraph into a single\n# run, then insert a new paragraph (\"Changes made by user of assignment 1\")\n# followed by a blank paragraph right after it.\n\n$d = $word.ActiveDocument\n\n$targetText = \">>>  your stuff after this line >>>\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = $targetText\n$find.Replacement.Text = $targetText + \"^pChanges made by user of assignment 1^p\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Forward = $true\n$find.Wrap = 1\n\n# wdReplaceAll = 2 -- rewrite the (single) match: this collapses the\n# paragraph's 3 runs into 1 run and adds the two new paragraphs right after.\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
